# "Arquivos da Semana da PÓS e Revisão da Proposta"
# Update the "CRONOGRAMA DISSERTAÇÃO" schedule: rework activity names/dates for
# rows 3-13, add two new activity rows (14-15), resize those rows, extend the
# used range to A1:D15, and move the active selection to the new last cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3 -----------------------------------------------------------------
$ws.Range("B3").Value = "Mapeamento  Sistemático da Literatura"
$ws.Range("C3").Value = "03/2016"
$ws.Range("D3").Value = "06/2016"
$ws.Rows.Item(3).RowHeight = 23.85

# --- Row 4 -----------------------------------------------------------------
$ws.Range("B4").Value = "Ponto de Controle 01 – Reunião com orientador sobre Revisão Sistemática da Literatura"
$ws.Range("C4").Value = "06/2016"
$ws.Range("D4").Value = "06/2016"
$ws.Rows.Item(4).RowHeight = 31.3

# --- Row 5 -----------------------------------------------------------------
$ws.Range("B5").Value = "Caracterização das Ferramentas de Gerenciamento de Requisição de Mudança"
$ws.Range("C5").Value = "06/2016"
$ws.Range("D5").Value = "07/2016"
$ws.Rows.Item(5).RowHeight = 31.3

# --- Row 6 -----------------------------------------------------------------
$ws.Range("B6").Value = "Ponto de Controle 02 – Reunião com orientador sobre a Caracterizão das FGRM"
$ws.Range("C6").Value = "07/2016"
$ws.Range("D6").Value = "07/2016"
$ws.Rows.Item(6).RowHeight = 31.3

# --- Row 7 -----------------------------------------------------------------
$ws.Range("B7").Value = "Pesquisa com Profissionais"
$ws.Range("C7").Value = "08/2016"
$ws.Range("D7").Value = "09/2016"
$ws.Rows.Item(7).RowHeight = 15.65

# --- Row 8 -----------------------------------------------------------------
$ws.Range("B8").Value = "Ponto de Controle 03 – Reunião com orientador sobre a Pesquisa com o Profissionais"
$ws.Range("C8").Value = "09/2016"
$ws.Range("D8").Value = "09/2016"
$ws.Rows.Item(8).RowHeight = 29.85

# --- Row 9 -----------------------------------------------------------------
$ws.Range("B9").Value = "Implementação da Ferramenta"
$ws.Range("C9").Value = "09/2016"
$ws.Range("D9").Value = "10/2016"
$ws.Rows.Item(9).RowHeight = 15.65

# --- Row 10 ----------------------------------------------------------------
$ws.Range("B10").Value = "Ponto de Controle 03 – Avaliação da Ferramenta Avaliada"
$ws.Range("C10").Value = "10/2016"
$ws.Range("D10").Value = "10/2016"
$ws.Rows.Item(10).RowHeight = 15.65

# --- Row 11 ----------------------------------------------------------------
$ws.Range("B11").Value = "Experimento de Avaliação da Ferramenta"
$ws.Range("C11").Value = "11/2016"
$ws.Range("D11").Value = "11/2016"
$ws.Rows.Item(11).RowHeight = 15.65

# --- Row 12 ----------------------------------------------------------------
$ws.Range("B12").Value = "Ponto de Controle 04 – Avaliação do Experimento junto com o orientador"
$ws.Range("C12").Value = "11/2016"
$ws.Range("D12").Value = "11/2016"
$ws.Rows.Item(12).RowHeight = 29.85

# --- Row 13 ----------------------------------------------------------------
$ws.Range("B13").Value = "Finalização do texto da dissertação"
$ws.Range("C13").Value = "12/2017"
$ws.Range("D13").Value = "12/2016"
$ws.Rows.Item(13).RowHeight = 15.65

# --- Row 14 (new) ------------------------------------------------------------
$ws.Range("A14").Value = 12
$ws.Range("B14").Value = "Ponto de Controle 05 – Avaliação do texto da dissertação com o orientador"
$ws.Range("C14").Value = "01/2017"
$ws.Range("D14").Value = "01/2017"
$ws.Range("A13:D13").Copy()
$ws.Range("A14:D14").PasteSpecial(-4122)
$ws.Rows.Item(14).RowHeight = 29.85

# --- Row 15 (new) ------------------------------------------------------------
$ws.Range("A15").Value = 13
$ws.Range("B15").Value = "Defesa da dissertação"
$ws.Range("C15").Value = "01/2017"
$ws.Range("D15").Value = "01/2017"
$ws.Range("A13:D13").Copy()
$ws.Range("A15:D15").PasteSpecial(-4122)
$ws.Rows.Item(15).RowHeight = 25.35

$excel.CutCopyMode = 0

# --- Selection ---------------------------------------------------------------
$ws.Range("B15").Select()
